# 1. Insert a new row at 66 (pushes existing rows 66+ down by one), so that a
#    second course ("Time Management") can be added right after the existing
#    "How to Meditate" course row (row 65), and fill in the previously empty
#    name/description cells for the "How to Meditate" row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(66).Insert()

# Fill in the previously-blank name/description for the existing course row 65
$ws.Range("B65").Value = "How to Meditate"
$ws.Range("C65").Value = "Learn how to meditate through experienced trainers."
# This course is now confirmed (is_confirmed 0 -> 1)
$ws.Range("O65").Value = 1

# New course row 66: "Time Management"
# (Excel's row-insert copies the formatting of the row above, so the
# venue_longitude/venue_latitude cells would otherwise inherit row 65's
# Tahoma-font style; reset them back to the default "Normal" style first.)
$ws.Range("E66:F66").Style = "Normal"
$ws.Range("A66").Value = 13
$ws.Range("B66").Value = "Time Management"
$ws.Range("C66").Value = "Learn how to manage your time with the help of proven time management experts."
$ws.Range("D66").Value = 4
$ws.Range("E66").Value = 14.56
$ws.Range("F66").Value = 120.99
$ws.Range("G66").Value = "2021-03-01"
$ws.Range("I66").Value = 75
$ws.Range("J66").Value = 1
$ws.Range("K66").Value = 5
$ws.Range("L66").Value = "Jake Trinity"
$ws.Range("M66").Value = 1
$ws.Range("N66").Value = 1
$ws.Range("O66").Value = 0
$ws.Range("P66").Value = 123

# 3. Add a new course_bookings row (352, course 13, client 301) for the newly
#    added "Time Management" course, appended after the existing booking row.
$ws.Range("A77").Value = 352
$ws.Range("B77").Value = 13
$ws.Range("C77").Value = 301

# 4. Update the "position" value in the ccpregistration sample row (row 55):
#    "cultural classes provider" -> "Cultural Class Provider"
$ws.Range("C55").Value = "Cultural Class Provider"

# 5. Leave the selection where the author last left it while entering this data.
$ws.Range("D67").Select()
